$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-23 with new property data
$ws.Range("A2").Value = "Al Majaz 2, Al Majaz, Sharjah"
$ws.Range("B2").Value = "48,000 AED/year"
$ws.Range("C2").Value = "2,500 sqft"
$ws.Range("D2").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-majaz-2-13153779.html"
$ws.Range("A3").Value = "Al Mahatta, Al Qasimia, Sharjah"
$ws.Range("B3").Value = "50,999 AED/year"
$ws.Range("C3").Value = "2,200 sqft"
$ws.Range("D3").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-qasimia-al-mahatta-12853729.html"
$ws.Range("A4").Value = "Al Mahatta, Al Qasimia, Sharjah"
$ws.Range("B4").Value = "54,999 AED/year"
$ws.Range("C4").Value = "3,500 sqft"
$ws.Range("D4").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-qasimia-al-mahatta-13133270.html"
$ws.Range("A5").Value = "Al Mahatta, Al Qasimia, Sharjah"
$ws.Range("B5").Value = "55,000 AED/year"
$ws.Range("C5").Value = "2,200 sqft"
$ws.Range("D5").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-qasimia-al-mahatta-12919247.html"
$ws.Range("A6").Value = "Al Wahda, Sharjah"
$ws.Range("B6").Value = "59,999 AED/year"
$ws.Range("C6").Value = "3,600 sqft"
$ws.Range("D6").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-wahda-13164090.html"
$ws.Range("A7").Value = "Al Majaz 3, Al Majaz, Sharjah"
$ws.Range("B7").Value = "79,999 AED/year"
$ws.Range("C7").Value = "3,000 sqft"
$ws.Range("D7").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-majaz-3-13013436.html"
$ws.Range("A8").Value = "Al Mamzar, Al Mamzar - Sharjah, Sharjah"
$ws.Range("B8").Value = "82,000 AED/year"
$ws.Range("C8").Value = "2,600 sqft"
$ws.Range("D8").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-mamzar-sharjah-al-mamzar-13027011.html"
$ws.Range("A9").Value = "Al Mamzar, Al Mamzar - Sharjah, Sharjah"
$ws.Range("B9").Value = "85,000 AED/year"
$ws.Range("C9").Value = "2,600 sqft"
$ws.Range("D9").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-mamzar-sharjah-al-mamzar-13027014.html"
$ws.Range("A10").Value = "Saeed Al Alami Building, Al Taawun, Sharjah"
$ws.Range("B10").Value = "90,000 AED/year"
$ws.Range("C10").Value = "3,000 sqft"
$ws.Range("D10").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-taawun-saeed-al-alami-building-12550137.html"
$ws.Range("A11").Value = "Sarh Al Emarat Tower, Al Majaz 3, Al Majaz, Sharjah"
$ws.Range("B11").Value = "90,000 AED/year"
$ws.Range("C11").Value = "2,034 sqft"
$ws.Range("D11").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-majaz-3-sarh-al-emarat-tower-13042633.html"
$ws.Range("A12").Value = "Majestic Tower, Al Taawun Street, Al Taawun, Sharjah"
$ws.Range("B12").Value = "94,999 AED/year"
$ws.Range("C12").Value = "3,000 sqft"
$ws.Range("D12").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-taawun-al-taawun-street-majestic-tower-12696062.html"
$ws.Range("A13").Value = "Al Maha Tower, Al Majaz, Sharjah"
$ws.Range("B13").Value = "95,000 AED/year"
$ws.Range("C13").Value = "3,000 sqft"
$ws.Range("D13").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-maha-tower-12822529.html"
$ws.Range("A14").Value = "Al Majaz 3, Al Majaz, Sharjah"
$ws.Range("B14").Value = "95,000 AED/year"
$ws.Range("C14").Value = "2,500 sqft"
$ws.Range("D14").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-majaz-3-10910376.html"
$ws.Range("A15").Value = "Al Mirgab, Al Heerah, Sharjah"
$ws.Range("B15").Value = "105,000 AED/year"
$ws.Range("C15").Value = "5,000 sqft"
$ws.Range("D15").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-heerah-al-mirgab-13150159.html"
$ws.Range("A16").Value = "Majestic Tower, Al Taawun Street, Al Taawun, Sharjah"
$ws.Range("B16").Value = "110,000 AED/year"
$ws.Range("C16").Value = "3,500 sqft"
$ws.Range("D16").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-taawun-al-taawun-street-majestic-tower-12762393.html"
$ws.Range("A17").Value = "Al Taawun Street, Al Taawun, Sharjah"
$ws.Range("B17").Value = "110,000 AED/year"
$ws.Range("C17").Value = "3,990 sqft"
$ws.Range("D17").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-taawun-al-taawun-street-12866617.html"
$ws.Range("A18").Value = "Al Majaz 3, Al Majaz, Sharjah"
$ws.Range("B18").Value = "115,000 AED/year"
$ws.Range("C18").Value = "3,200 sqft"
$ws.Range("D18").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-majaz-3-11998377.html"
$ws.Range("A19").Value = "Al Khan Corniche, Al Khan, Sharjah"
$ws.Range("B19").Value = "120,000 AED/year"
$ws.Range("C19").Value = "3,510 sqft"
$ws.Range("D19").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-khan-al-khan-corniche-12803583.html"
$ws.Range("A20").Value = "Al Majaz 3, Al Majaz, Sharjah"
$ws.Range("B20").Value = "120,000 AED/year"
$ws.Range("C20").Value = "3,500 sqft"
$ws.Range("D20").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-majaz-3-12500554.html"
$ws.Range("A21").Value = "Al Majaz 3, Al Majaz, Sharjah"
$ws.Range("B21").Value = "120,000 AED/year"
$ws.Range("C21").Value = "3,555 sqft"
$ws.Range("D21").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-majaz-3-13149354.html"
$ws.Range("A22").Value = "Al Maha Tower, Al Majaz, Sharjah"
$ws.Range("B22").Value = "125,000 AED/year"
$ws.Range("C22").Value = "3,000 sqft"
$ws.Range("D22").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-sharjah-al-majaz-al-maha-tower-11951702.html"
$ws.Range("A23").Value = "Greenview, EMAAR South, Dubai South (Dubai World Central), Dubai"
$ws.Range("B23").Value = "140,000 AED/year"
$ws.Range("C23").Value = "2,540 sqft"
$ws.Range("D23").Value = "https://www.propertyfinder.ae/en/plp/rent/apartment-for-rent-dubai-dubai-south-dubai-world-central-emaar-south-greenview-12916545.html"

# Remove row 24 (no longer present in the updated dataset); this shifts
# nothing below it (it was the last row) and shrinks the used range to A1:D23
$ws.Rows.Item(24).Delete()
